# Renames the inline-picture "name" (wp:docPr/@name and pic:cNvPr/@name)
# for the three logo images in the headers/footers of the document:
#   - BTec logo (header, first page):      image2.jpg -> image1.jpg
#   - Pearson logo (footer, default page): image1.png -> image2.png
#   - Pearson logo (footer, first page):   image1.png -> image2.png
#
# The InlineShape object model only exposes a `.Name` setter that updates
# wp:docPr/@name, it does not keep pic:cNvPr/@name synchronised. To update
# both attributes consistently (as a real Word save would), we round-trip
# the shape's own Range through WordOpenXML, patch the "name=" attribute
# values for that shape only, and feed the patched XML back in with
# InsertXML (this replaces the shape range in place).

function Rename-InlineShape {
    param(
        $shape,
        [string]$oldName,
        [string]$newName
    )

    $rng = $shape.Range
    $xml = $rng.WordOpenXML

    $oldAttr = 'name="' + $oldName + '"'
    $newAttr = 'name="' + $newName + '"'

    $patched = $xml.Replace($oldAttr, $newAttr)
    $rng.InsertXML($patched)
}

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# BTec logo lives in the "first page" header (Header.Item(2) in this doc).
$btecHeader = $sec.Headers.Item(2)
$btecShape = $btecHeader.Range.InlineShapes.Item(1)
Rename-InlineShape -shape $btecShape -oldName "image2.jpg" -newName "image1.jpg"

# Pearson logo in the "default" footer (Footer.Item(1) in this doc).
$pearsonFooterDefault = $sec.Footers.Item(1)
$pearsonShapeDefault = $pearsonFooterDefault.Range.InlineShapes.Item(1)
Rename-InlineShape -shape $pearsonShapeDefault -oldName "image1.png" -newName "image2.png"

# Pearson logo in the "first page" footer (Footer.Item(2) in this doc).
$pearsonFooterFirst = $sec.Footers.Item(2)
$pearsonShapeFirst = $pearsonFooterFirst.Range.InlineShapes.Item(1)
Rename-InlineShape -shape $pearsonShapeFirst -oldName "image1.png" -newName "image2.png"
